$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.8606349999999999
$ws.Cells.Item(2, 8).Value = 2.581905
$ws.Cells.Item(2, 9).Value = 0.0262626340301864
$ws.Cells.Item(2, 10).Value = 0.0262626340301864
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.10121233333333
$ws.Cells.Item(2, 14).Value = 39.303637
$ws.Cells.Item(2, 15).Value = 0.1081423012186565
$ws.Cells.Item(2, 16).Value = 0.1081423012186565
$ws.Cells.Item(2, 17).Value = 11.27536187649833
$ws.Cells.Item(2, 18).Value = 101.478256888485
$ws.Cells.Item(2, 19).Value = 0.002840101680087755
$ws.Cells.Item(2, 20).Value = 0.002840101680087755
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.8606349999999999
$ws.Cells.Item(3, 8).Value = 2.581905
$ws.Cells.Item(3, 9).Value = 0.0262626340301864
$ws.Cells.Item(3, 10).Value = 0.0262626340301864
$ws.Cells.Item(3, 15).Value = 0.5751439322003361
$ws.Cells.Item(3, 16).Value = 0.5751439322003362
$ws.Cells.Item(3, 17).Value = 59.96687599165165
$ws.Cells.Item(3, 18).Value = 539.7018839248649
$ws.Cells.Item(3, 19).Value = 0.01510479460605977
$ws.Cells.Item(3, 20).Value = 0.01510479460605977
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.8606349999999999
$ws.Cells.Item(4, 8).Value = 2.581905
$ws.Cells.Item(4, 9).Value = 0.0262626340301864
$ws.Cells.Item(4, 10).Value = 0.0262626340301864
$ws.Cells.Item(4, 13).Value = 38.36920666666666
$ws.Cells.Item(4, 14).Value = 115.10762
$ws.Cells.Item(4, 15).Value = 0.3167137665810073
$ws.Cells.Item(4, 16).Value = 0.3167137665810074
$ws.Cells.Item(4, 17).Value = 33.02188217956666
$ws.Cells.Item(4, 18).Value = 297.1969396161
$ws.Cells.Item(4, 19).Value = 0.008317737744038874
$ws.Cells.Item(4, 20).Value = 0.008317737744038876
$ws.Cells.Item(5, 8).Value = 58.40949000000001
$ws.Cells.Item(5, 9).Value = 0.5941299388474139
$ws.Cells.Item(5, 10).Value = 0.5941299388474139
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 13.10121233333333
$ws.Cells.Item(5, 14).Value = 39.303637
$ws.Cells.Item(5, 15).Value = 0.1081423012186565
$ws.Cells.Item(5, 16).Value = 0.1081423012186565
$ws.Cells.Item(5, 17).Value = 255.0783769239034
$ws.Cells.Item(5, 18).Value = 2295.70539231513
$ws.Cells.Item(5, 19).Value = 0.06425057880985899
$ws.Cells.Item(5, 20).Value = 0.06425057880985899
$ws.Cells.Item(6, 8).Value = 58.40949000000001
$ws.Cells.Item(6, 9).Value = 0.5941299388474139
$ws.Cells.Item(6, 10).Value = 0.5941299388474139
$ws.Cells.Item(6, 15).Value = 0.5751439322003361
$ws.Cells.Item(6, 16).Value = 0.5751439322003362
$ws.Cells.Item(6, 19).Value = 0.3417102292666469
$ws.Cells.Item(6, 20).Value = 0.341710229266647
$ws.Cells.Item(7, 8).Value = 58.40949000000001
$ws.Cells.Item(7, 9).Value = 0.5941299388474139
$ws.Cells.Item(7, 10).Value = 0.5941299388474139
$ws.Cells.Item(7, 13).Value = 38.36920666666666
$ws.Cells.Item(7, 14).Value = 115.10762
$ws.Cells.Item(7, 15).Value = 0.3167137665810073
$ws.Cells.Item(7, 16).Value = 0.3167137665810074
$ws.Cells.Item(7, 17).Value = 747.0419310348667
$ws.Cells.Item(7, 18).Value = 6723.377379313801
$ws.Cells.Item(7, 19).Value = 0.188169130770908
$ws.Cells.Item(7, 20).Value = 0.1881691307709081
$ws.Cells.Item(8, 7).Value = 12.439858
$ws.Cells.Item(8, 8).Value = 37.319574
$ws.Cells.Item(8, 9).Value = 0.3796074271223998
$ws.Cells.Item(8, 10).Value = 0.3796074271223997
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 13.10121233333333
$ws.Cells.Item(8, 14).Value = 39.303637
$ws.Cells.Item(8, 15).Value = 0.1081423012186565
$ws.Cells.Item(8, 16).Value = 0.1081423012186565
$ws.Cells.Item(8, 17).Value = 162.9772210545153
$ws.Cells.Item(8, 18).Value = 1466.794989490638
$ws.Cells.Item(8, 19).Value = 0.04105162072870974
$ws.Cells.Item(8, 20).Value = 0.04105162072870974
$ws.Cells.Item(9, 7).Value = 12.439858
$ws.Cells.Item(9, 8).Value = 37.319574
$ws.Cells.Item(9, 9).Value = 0.3796074271223998
$ws.Cells.Item(9, 10).Value = 0.3796074271223997
$ws.Cells.Item(9, 15).Value = 0.5751439322003361
$ws.Cells.Item(9, 16).Value = 0.5751439322003362
$ws.Cells.Item(9, 17).Value = 866.7779279715046
$ws.Cells.Item(9, 18).Value = 7801.00135174354
$ws.Cells.Item(9, 19).Value = 0.2183289083276295
$ws.Cells.Item(9, 20).Value = 0.2183289083276295
$ws.Cells.Item(10, 7).Value = 12.439858
$ws.Cells.Item(10, 8).Value = 37.319574
$ws.Cells.Item(10, 9).Value = 0.3796074271223998
$ws.Cells.Item(10, 10).Value = 0.3796074271223997
$ws.Cells.Item(10, 13).Value = 38.36920666666666
$ws.Cells.Item(10, 14).Value = 115.10762
$ws.Cells.Item(10, 15).Value = 0.3167137665810073
$ws.Cells.Item(10, 16).Value = 0.3167137665810074
$ws.Cells.Item(10, 17).Value = 477.3074825059866
$ws.Cells.Item(10, 18).Value = 4295.76734255388
$ws.Cells.Item(10, 19).Value = 0.1202268980660605
$ws.Cells.Item(10, 20).Value = 0.1202268980660605
